# The deck currently carries the custom "Integral" theme on its slide
# master/design. This restores the design to the stock default "Office
# Theme" colour palette.
#
# PowerPoint's object model doesn't expose a way to rewrite a whole theme
# part in one shot - the supported mechanism is to edit the live theme's
# colour scheme slot-by-slot via ThemeColorScheme.Colors(i).RGB (fonts would
# go through ThemeFontScheme.MajorFont/MinorFont, left untouched here since
# both themes already share the same "Office" font scheme).

$p   = $ppt.ActivePresentation
$sm  = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

function Set-ThemeColor($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    # VBA RGB() packs as R + G*256 + B*65536
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Office Theme default colour scheme, in MSO theme-colour-index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
Set-ThemeColor $tcs 1  "000000"
Set-ThemeColor $tcs 2  "FFFFFF"
Set-ThemeColor $tcs 3  "44546A"
Set-ThemeColor $tcs 4  "E7E6E6"
Set-ThemeColor $tcs 5  "5B9BD5"
Set-ThemeColor $tcs 6  "ED7D31"
Set-ThemeColor $tcs 7  "A5A5A5"
Set-ThemeColor $tcs 8  "FFC000"
Set-ThemeColor $tcs 9  "4472C4"
Set-ThemeColor $tcs 10 "70AD47"
Set-ThemeColor $tcs 11 "0563C1"
Set-ThemeColor $tcs 12 "954F72"
